$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry describes the new values for a row (B=Coin, C=Link, D=Price, E=Volume(1h)).
# A brand-new coin (OKB) is inserted at row 8; every row from the old row 8 (Cardano)
# through the old row 50 (TrueUSD) shifts down by one, and the old last row (USDD)
# drops off the bottom of the list. D holds price text that looks numeric (e.g.
# "24.856.95", "0.9979", "1.000"), so it must be written as Text to avoid Excel
# coercing it into a number and losing the original formatting/precision.
$updates = @(
    @{Row=2; D='24.856.95'; E='  -4.27%  '},
    @{Row=3; D='1.634.18'; E='  -6.42%  '},
    @{Row=4; D='0.9979'; E='  -0.18%  '},
    @{Row=5; D='234.55'; E='  -5.69%  '},
    @{Row=6; D='1.000'; E='  +0.03%  '},
    @{Row=7; D='0.4707'; E='  -6.66%  '},
    @{Row=8; B='OKB'; C='https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D='38.77'; E='  -4.60%  '},
    @{Row=9; B='Cardano'; C='https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'; D='0.2554'; E='  -6.92%  '},
    @{Row=10; B='Dogecoin'; C='https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; D='0.06059'; E='  -2.03%  '},
    @{Row=11; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.06946'; E='  -4.54%  '},
    @{Row=12; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='1.637.92'; E='  -6.11%  '},
    @{Row=13; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='14.58'; E='  -4.10%  '},
    @{Row=14; B='Polygon'; C='https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; D='0.6044'; E='  -7.58%  '},
    @{Row=15; B='Polkadot'; C='https://coinranking.com/coin/25W7FG7om+polkadot-dot'; D='4.326'; E='  -6.96%  '},
    @{Row=16; B='Litecoin'; C='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; D='72.75'; E='  -6.39%  '},
    @{Row=17; B='Dai'; C='https://coinranking.com/coin/MoTuySvg7+dai-dai'; D='1.000'; E='  +0.01%  '},
    @{Row=18; B='BinanceUSD'; C='https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; D='0.9992'; E='  -0.03%  '},
    @{Row=19; B='WrappedBTC'; C='https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; D='24.868.05'; E='  -4.29%  '},
    @{Row=20; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.000006563'; E='  -4.00%  '},
    @{Row=21; B='Avalanche'; C='https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; D='11.10'; E='  -6.25%  '},
    @{Row=22; B='WrappedliquidstakedEther2.0'; C='https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; D='1.847.96'; E='  -6.07%  '},
    @{Row=23; B='Uniswap'; C='https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; D='4.347'; E='  -1.73%  '},
    @{Row=24; B='Cosmos'; C='https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D='8.551'; E='  -1.97%  '},
    @{Row=25; B='Chainlink'; C='https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; D='5.207'; E='  -3.54%  '},
    @{Row=26; B='Monero'; C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D='133.11'; E='  -2.71%  '},
    @{Row=27; B='EthereumClassic'; C='https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D='14.73'; E='  -3.32%  '},
    @{Row=28; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='1.380'; E='  -8.29%  '},
    @{Row=29; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='103.21'; E='  -2.21%  '},
    @{Row=30; B='LidoDAOToken'; C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; D='1.625'; E='  -8.70%  '},
    @{Row=31; B='InternetComputer(DFINITY)'; C='https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; D='3.772'; E='  -2.31%  '},
    @{Row=32; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.07715'; E='  -6.02%  '},
    @{Row=33; B='Filecoin'; C='https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D='3.523'; E='  -3.08%  '},
    @{Row=34; B='Frax'; C='https://coinranking.com/coin/KfWtaeV1W+frax-frax'; D='0.9990'; E='  -0.01%  '},
    @{Row=35; B='Hedera'; C='https://coinranking.com/coin/jad286TjB+hedera-hbar'; D='0.04288'; E='  -8.24%  '},
    @{Row=36; B='HuobiToken'; C='https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; D='2.582'; E='  -2.77%  '},
    @{Row=37; B='ARBITRUM'; C='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D='0.9183'; E='  -7.61%  '},
    @{Row=38; B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='0.5783'; E='  -6.71%  '},
    @{Row=39; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='2.533'; E='  -7.91%  '},
    @{Row=40; B='VeChain'; C='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; D='0.01535'; E='  -4.78%  '},
    @{Row=41; B='PaxDollar'; C='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'; D='0.9988'; E='  -0.11%  '},
    @{Row=42; B='TrustWalletToken'; C='https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; D='0.8145'; E='  +7.20%  '},
    @{Row=43; B='Quant'; C='https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; D='97.02'; E='  -3.51%  '},
    @{Row=44; B='RenderToken'; C='https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; D='1.769'; E='  -8.18%  '},
    @{Row=45; B='TheSandbox'; C='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; D='0.3686'; E='  -6.44%  '},
    @{Row=46; B='FraxShare'; C='https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; D='4.699'; E='  -6.16%  '},
    @{Row=47; D='0.05201'; E='  -1.40%  '},
    @{Row=48; B='Algorand'; C='https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; D='0.1087'; E='  -5.45%  '},
    @{Row=49; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='6.007'; E='  -4.80%  '},
    @{Row=50; B='Elrond'; C='https://coinranking.com/coin/omwkOTglq+elrond-egld'; D='29.37'; E='  -4.20%  '},
    @{Row=51; B='TrueUSD'; C='https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'; D='1.000'; E='  -0.21%  '}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Cells.Item($r, 2).Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Cells.Item($r, 3).Value = $u.C }
    if ($u.ContainsKey('D')) {
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value = $u.E }
}
